# Apply the "Add, Update, Delete category and product feature" change:
# Insert a new "Cost" column (= Price - Price*20%) into the Product sheet
# between Price and Quantity, and update the view/selection state.

$wb = $excel.ActiveWorkbook
$wsCategory = $wb.Worksheets.Item("Category")
$wsProduct  = $wb.Worksheets.Item("Product")

# --- Product sheet: insert a new column E ("Cost") --------------------
# This shifts the existing Quantity/ImageName columns from E/F to F/G.
$wsProduct.Columns.Item(5).Insert() | Out-Null

# Header
$wsProduct.Range("E1").Value = "Cost"

# Cost formula for every data row (2-65): Price minus 20%.
$wsProduct.Range("E2").Formula = "=D2 - (D2 * 20 / 100)"
$wsProduct.Range("E3:E65").Formula = "=D3 - (D3 * 20 / 100)"

# --- View / selection state --------------------------------------------
# Category sheet keeps its previous selection (D33) but is no longer the
# active tab / no longer scrolled to the top.
$wsCategory.Range("D33").Select() | Out-Null

# Product sheet becomes the active tab with a new selection (M34).
$wsProduct.Range("M34").Select() | Out-Null
$wsProduct.Activate() | Out-Null
